# update sfserver update experiments
# Fill in the "aggregation" latency column (C) on the "eventdriven" sheet
# for each task-number block (2, 4, 8, 16, 32, 64 tasks), mirroring the
# per-run samples already recorded in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eventdriven")
$ws.Activate()

# task number = 2 block
$ws.Range("C3").Value = 2.1328719999999999
$ws.Range("C4").Value = 2.1349960000000001
$ws.Range("C5").Value = 2.1153960000000001
$ws.Range("C6").Value = 2.182798
$ws.Range("C7").Value = 2.1255259999999998

# task number = 4 block
$ws.Range("C10").Value = 4
$ws.Range("C11").Value = 2.1342910000000002
$ws.Range("C12").Value = 2.1290909999999998
$ws.Range("C13").Value = 2.1300340000000002
$ws.Range("C14").Value = 2.12243
$ws.Range("C15").Value = 2.1209600000000002

# task number = 8 block
$ws.Range("C19").Value = 8
$ws.Range("C20").Value = 2.140355
$ws.Range("C21").Value = 2.1380560000000002
$ws.Range("C22").Value = 2.139481
$ws.Range("C23").Value = 2.1326369999999999
$ws.Range("C24").Value = 2.1319370000000002

# task number = 16 block
$ws.Range("C27").Value = 16
$ws.Range("C28").Value = 2.1677909999999998
$ws.Range("C29").Value = 2.1599439999999999
$ws.Range("C30").Value = 2.1534990000000001
$ws.Range("C31").Value = 2.149159
$ws.Range("C32").Value = 2.1557819999999999

# task number = 32 block
$ws.Range("C36").Value = 32
$ws.Range("C37").Value = 2.2002830000000002
$ws.Range("C38").Value = 2.1922450000000002
$ws.Range("C39").Value = 2.1933370000000001
$ws.Range("C40").Value = 2.1951710000000002
$ws.Range("C41").Value = 2.1929270000000001

# task number = 64 block
$ws.Range("C46").Value = 64

# Leave the view scrolled down to the newly-added data, matching where
# the author was working when the sheet was saved.
$excel.Goto($ws.Range("A14"), $true) | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("C46").Select() | Out-Null
